# Auto-generated edit script to apply scheduled-runner data refresh
# to the Typhon_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 40004384  # H113
$ws.Cells.Item(113, 9).Value = 90912780  # I113
$ws.Cells.Item(113, 10).Value = 4928.5713  # J113
$ws.Cells.Item(113, 11).Value = 90912780  # K113
$ws.Cells.Item(113, 12).Value = 4928.5713  # L113
$ws.Cells.Item(113, 13).Value = -90909526  # M113
$ws.Cells.Item(113, 14).Value = -11436.5713  # N113
$ws.Cells.Item(129, 8).Value = 371546.56  # H129
$ws.Cells.Item(129, 9).Value = 336.7143  # I129
$ws.Cells.Item(129, 10).Value = 501470  # J129
$ws.Cells.Item(129, 11).Value = 1010.1429  # K129
$ws.Cells.Item(129, 12).Value = 1504410  # L129
$ws.Cells.Item(129, 13).Value = 3989.8571  # M129
$ws.Cells.Item(129, 14).Value = -1514410  # N129
$ws.Cells.Item(132, 8).Value = 1943.234  # H132
$ws.Cells.Item(132, 9).Value = 2146.5  # I132
$ws.Cells.Item(132, 10).Value = 781.7143  # J132
$ws.Cells.Item(132, 11).Value = 6439.5  # K132
$ws.Cells.Item(132, 12).Value = 2345.1429  # L132
$ws.Cells.Item(132, 13).Value = -3909.5  # M132
$ws.Cells.Item(132, 14).Value = -7405.1429  # N132
$ws.Cells.Item(137, 8).Value = 16293.418  # H137
$ws.Cells.Item(137, 9).Value = 1209.5834  # I137
$ws.Cells.Item(137, 10).Value = 54399.95  # J137
$ws.Cells.Item(137, 11).Value = 3628.7502  # K137
$ws.Cells.Item(137, 12).Value = 163199.85  # L137
$ws.Cells.Item(137, 13).Value = -1078.7502  # M137
$ws.Cells.Item(137, 14).Value = -168299.85  # N137
$ws.Cells.Item(138, 8).Value = 12823038  # H138
$ws.Cells.Item(138, 10).Value = 2468.0144  # J138
$ws.Cells.Item(138, 12).Value = 7404.0432  # L138
$ws.Cells.Item(138, 14).Value = -17684.0432  # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 54628.305  # H32
$ws.Cells.Item(32, 9).Value = 62221.9  # I32
$ws.Cells.Item(32, 11).Value = 62221.9  # K32
$ws.Cells.Item(32, 13).Value = -61934.9  # M32
$ws.Cells.Item(45, 8).Value = 2979.2144  # H45
$ws.Cells.Item(45, 9).Value = 2747.1428  # I45
$ws.Cells.Item(45, 10).Value = 3675.4285  # J45
$ws.Cells.Item(45, 11).Value = 2747.1428  # K45
$ws.Cells.Item(45, 12).Value = 3675.4285  # L45
$ws.Cells.Item(45, 13).Value = -2370.1428  # M45
$ws.Cells.Item(45, 14).Value = -4429.4285  # N45
$ws.Cells.Item(61, 8).Value = 1374.45  # H61
$ws.Cells.Item(61, 9).Value = 1085.5883  # I61
$ws.Cells.Item(61, 10).Value = 3011.3333  # J61
$ws.Cells.Item(61, 11).Value = 1085.5883  # K61
$ws.Cells.Item(61, 12).Value = 3011.3333  # L61
$ws.Cells.Item(61, 13).Value = -873.5882999999999  # M61
$ws.Cells.Item(61, 14).Value = -3435.3333  # N61
$ws.Cells.Item(97, 8).Value = 3035  # H97
$ws.Cells.Item(97, 9).Value = 3035  # I97
$ws.Cells.Item(97, 10).Value = 0  # J97
$ws.Cells.Item(97, 11).Value = 3035  # K97
$ws.Cells.Item(97, 12).Value = 0  # L97
$ws.Cells.Item(97, 13).Value = -2539  # M97
$ws.Cells.Item(97, 14).ClearContents()  # N97
$ws.Cells.Item(102, 8).Value = 1633.3334  # H102
$ws.Cells.Item(102, 9).Value = 900  # I102
$ws.Cells.Item(102, 11).Value = 900  # K102
$ws.Cells.Item(102, 13).Value = 722  # M102
$ws.Cells.Item(122, 8).Value = 1985.2069  # H122
$ws.Cells.Item(122, 9).Value = 1708.8125  # I122
$ws.Cells.Item(122, 11).Value = 5126.4375  # K122
$ws.Cells.Item(122, 13).Value = -2676.4375  # M122
$ws.Cells.Item(136, 8).Value = 1374.45  # H136
$ws.Cells.Item(136, 9).Value = 1085.5883  # I136
$ws.Cells.Item(136, 10).Value = 3011.3333  # J136
$ws.Cells.Item(136, 11).Value = 3256.7649  # K136
$ws.Cells.Item(136, 12).Value = 9033.999899999999  # L136
$ws.Cells.Item(136, 13).Value = -706.7648999999997  # M136
$ws.Cells.Item(136, 14).Value = -14133.9999  # N136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(38, 8).Value = 0  # H38
$ws.Cells.Item(38, 10).Value = 0  # J38
$ws.Cells.Item(38, 12).Value = 0  # L38
$ws.Cells.Item(38, 14).ClearContents()  # N38
$ws.Cells.Item(99, 8).Value = 1431.0667  # H99
$ws.Cells.Item(99, 9).Value = 1618.3334  # I99
$ws.Cells.Item(99, 11).Value = 1618.3334  # K99
$ws.Cells.Item(99, 13).Value = -120.3334  # M99

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 11231  # H132
$ws.Cells.Item(132, 9).Value = 12866.884  # I132
$ws.Cells.Item(132, 10).Value = 4196.7  # J132
$ws.Cells.Item(132, 11).Value = 38600.652  # K132
$ws.Cells.Item(132, 12).Value = 12590.1  # L132
$ws.Cells.Item(132, 13).Value = -36070.652  # M132
$ws.Cells.Item(132, 14).Value = -17650.1  # N132
$ws.Cells.Item(134, 8).Value = 844.1622  # H134
$ws.Cells.Item(134, 9).Value = 677.8  # I134
$ws.Cells.Item(134, 11).Value = 2033.4  # K134
$ws.Cells.Item(134, 13).Value = 501.6000000000001  # M134

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 3183.0833  # H56
$ws.Cells.Item(56, 9).Value = 3183.0833  # I56
$ws.Cells.Item(56, 11).Value = 3183.0833  # K56
$ws.Cells.Item(56, 13).Value = -2653.0833  # M56
$ws.Cells.Item(117, 8).Value = 1822.375  # H117
$ws.Cells.Item(117, 9).Value = 739.5  # I117
$ws.Cells.Item(117, 10).Value = 2905.25  # J117
$ws.Cells.Item(117, 11).Value = 2218.5  # K117
$ws.Cells.Item(117, 12).Value = 8715.75  # L117
$ws.Cells.Item(117, 13).Value = 1223.5  # M117
$ws.Cells.Item(117, 14).Value = -15599.75  # N117
$ws.Cells.Item(118, 8).Value = 55561100  # H118
$ws.Cells.Item(118, 9).Value = 125000220  # I118
$ws.Cells.Item(118, 10).Value = 9799.799999999999  # J118
$ws.Cells.Item(118, 11).Value = 375000660  # K118
$ws.Cells.Item(118, 12).Value = 29399.4  # L118
$ws.Cells.Item(118, 13).Value = -374999417  # M118
$ws.Cells.Item(118, 14).Value = -31885.4  # N118
$ws.Cells.Item(131, 8).Value = 750.63  # H131
$ws.Cells.Item(131, 9).Value = 313  # I131
$ws.Cells.Item(131, 10).Value = 799.25555  # J131
$ws.Cells.Item(131, 11).Value = 939  # K131
$ws.Cells.Item(131, 12).Value = 2397.76665  # L131
$ws.Cells.Item(131, 13).Value = 4101  # M131
$ws.Cells.Item(131, 14).Value = -12477.76665  # N131

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 7980  # H33
$ws.Cells.Item(33, 10).Value = 7980  # J33
$ws.Cells.Item(33, 12).Value = 7980  # L33
$ws.Cells.Item(33, 14).Value = -8484  # N33
$ws.Cells.Item(97, 8).Value = 990.2857  # H97
$ws.Cells.Item(97, 10).Value = 1250  # J97
$ws.Cells.Item(97, 12).Value = 1250  # L97
$ws.Cells.Item(97, 14).Value = -2242  # N97
$ws.Cells.Item(107, 8).Value = 5128431.5  # H107
$ws.Cells.Item(107, 9).Value = 300.33334  # I107
$ws.Cells.Item(107, 10).Value = 12820628  # J107
$ws.Cells.Item(107, 11).Value = 300.33334  # K107
$ws.Cells.Item(107, 12).Value = 12820628  # L107
$ws.Cells.Item(107, 13).Value = 1619.66666  # M107
$ws.Cells.Item(107, 14).Value = -12824468  # N107
$ws.Cells.Item(121, 8).Value = 30305  # H121
$ws.Cells.Item(121, 10).Value = 30305  # J121
$ws.Cells.Item(121, 12).Value = 30305  # L121
$ws.Cells.Item(121, 14).Value = -33799  # N121
$ws.Cells.Item(122, 8).Value = 95238860  # H122
$ws.Cells.Item(122, 9).Value = 27778626  # I122
$ws.Cells.Item(122, 10).Value = 500000300  # J122
$ws.Cells.Item(122, 11).Value = 83335878  # K122
$ws.Cells.Item(122, 12).Value = 1500000900  # L122
$ws.Cells.Item(122, 13).Value = -83333428  # M122
$ws.Cells.Item(122, 14).Value = -1500005800  # N122

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 893.41174  # H93
$ws.Cells.Item(93, 9).Value = 791.38464  # I93
$ws.Cells.Item(93, 10).Value = 1225  # J93
$ws.Cells.Item(93, 11).Value = 791.38464  # K93
$ws.Cells.Item(93, 12).Value = 1225  # L93
$ws.Cells.Item(93, 13).Value = 456.61536  # M93
$ws.Cells.Item(93, 14).Value = -3721  # N93
$ws.Cells.Item(140, 8).Value = 50009.2  # H140
$ws.Cells.Item(140, 10).Value = 50009.2  # J140
$ws.Cells.Item(140, 12).Value = 50009.2  # L140
$ws.Cells.Item(140, 14).Value = -60369.2  # N140

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1291.4706  # H122
$ws.Cells.Item(122, 9).Value = 1126  # I122
$ws.Cells.Item(122, 10).Value = 1527.8572  # J122
$ws.Cells.Item(122, 11).Value = 3378  # K122
$ws.Cells.Item(122, 12).Value = 4583.571599999999  # L122
$ws.Cells.Item(122, 13).Value = -928  # M122
$ws.Cells.Item(122, 14).Value = -9483.571599999999  # N122

